# AI can build town on captured hexa
# Updates the "AI" worksheet: renames a hexa row to "Posvatna hora",
# retitles the duplicate "Points" header columns (O1/AZ1) to "Capture",
# rewrites several numeric stats for row 2, and moves the active
# selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AI")

# Row 2 key / label: numeric 8 -> text "Posvatna hora" (A2)
$ws.Range("A2").Value = "Posvatna hora"

# Column headers (row 1): "Points" -> "Capture" for the capture-hexa group
$ws.Range("O1").Value = "Capture"
$ws.Range("AZ1").Value = "Capture"

# Row 2 numeric updates
$ws.Range("B2").Value = 100
$ws.Range("H2").Value = 20000
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 90
$ws.Range("P2").Value = 5

$ws.Range("AS2").Value = 1000
$ws.Range("AW2").Value = 5000
$ws.Range("AX2").Value = 10
$ws.Range("AY2").Value = 0
$ws.Range("AZ2").Value = 90
$ws.Range("BA2").Value = 0

$ws.Range("BT2").Value = 1

# Move the active selection to C3
$ws.Range("C3").Select() | Out-Null
